# Daily attendance processing - 2026-01-29 20:44:37
#
# The "Recorded By" column (G) lists the users who recorded/edited a
# session's attendance. For every session row where it currently reads
# "System, dnasr281@gmail.com", the order of the two names is swapped to
# "dnasr281@gmail.com, System".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldText = "System, dnasr281@gmail.com"
$newText = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # Column G = "Recorded By"
    $value = $cell.Value()
    if ($value -eq $oldText) {
        $cell.Value = $newText
    }
}
